$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5860.6
$ws.Range("I64").Value = 5266.6665
$ws.Range("K64").Value = 5266.6665
$ws.Range("M64").Value = -5018.6665
$ws.Range("H67").Value = 5860.6
$ws.Range("I67").Value = 5266.6665
$ws.Range("K67").Value = 5266.6665
$ws.Range("M67").Value = -4408.6665
$ws.Range("H74").Value = 12000
$ws.Range("I74").Value = 15000
$ws.Range("J74").Value = 9000
$ws.Range("K74").Value = 15000
$ws.Range("L74").Value = 9000
$ws.Range("M74").Value = -14064
$ws.Range("N74").Value = -10872
$ws.Range("H76").Value = 111113750
$ws.Range("I76").Value = 142859710
$ws.Range("K76").Value = 142859710
$ws.Range("M76").Value = -142859395
$ws.Range("H77").Value = 12000
$ws.Range("I77").Value = 15000
$ws.Range("J77").Value = 9000
$ws.Range("K77").Value = 75000
$ws.Range("L77").Value = 45000
$ws.Range("M77").Value = -70320
$ws.Range("N77").Value = -54360
$ws.Range("H79").Value = 111113750
$ws.Range("I79").Value = 142859710
$ws.Range("K79").Value = 142859710
$ws.Range("M79").Value = -142858618
$ws.Range("H132").Value = 1685169
$ws.Range("I132").Value = 1419.3846
$ws.Range("J132").Value = 7939096
$ws.Range("K132").Value = 4258.1538
$ws.Range("L132").Value = 23817288
$ws.Range("M132").Value = -1728.1538
$ws.Range("N132").Value = -23822348
$ws.Range("H137").Value = 18763404
$ws.Range("I137").Value = 3907438.5
$ws.Range("J137").Value = 78187260
$ws.Range("K137").Value = 11722315.5
$ws.Range("L137").Value = 234561780
$ws.Range("M137").Value = -11719765.5
$ws.Range("N137").Value = -234566880
$ws.Range("H138").Value = 1978.8506
$ws.Range("I138").Value = 1261.0172
$ws.Range("J138").Value = 3414.5173
$ws.Range("K138").Value = 3783.0516
$ws.Range("L138").Value = 10243.5519
$ws.Range("M138").Value = 1356.9484
$ws.Range("N138").Value = -20523.5519
$ws.Range("H141").Value = 1484.3214
$ws.Range("I141").Value = 1006.5
$ws.Range("J141").Value = 4351.25
$ws.Range("K141").Value = 3019.5
$ws.Range("L141").Value = 13053.75
$ws.Range("M141").Value = 2160.5
$ws.Range("N141").Value = -23413.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1001133.4
$ws.Range("I45").Value = 1251065.4
$ws.Range("J45").Value = 1405.5
$ws.Range("K45").Value = 1251065.4
$ws.Range("L45").Value = 1405.5
$ws.Range("M45").Value = -1250688.4
$ws.Range("N45").Value = -2159.5
$ws.Range("H61").Value = 3351164.8
$ws.Range("I61").Value = 1667852.6
$ws.Range("J61").Value = 11767726
$ws.Range("K61").Value = 1667852.6
$ws.Range("L61").Value = 11767726
$ws.Range("M61").Value = -1667640.6
$ws.Range("N61").Value = -11768150
$ws.Range("H74").Value = 34849870
$ws.Range("I74").Value = 33334110
$ws.Range("J74").Value = 38097932
$ws.Range("K74").Value = 33334110
$ws.Range("L74").Value = 38097932
$ws.Range("M74").Value = -33333236
$ws.Range("N74").Value = -38099680
$ws.Range("H77").Value = 34849870
$ws.Range("I77").Value = 33334110
$ws.Range("J77").Value = 38097932
$ws.Range("K77").Value = 166670550
$ws.Range("L77").Value = 190489660
$ws.Range("M77").Value = -166666182
$ws.Range("N77").Value = -190498396
$ws.Range("H132").Value = 9262725
$ws.Range("I132").Value = 11115286
$ws.Range("K132").Value = 33345858
$ws.Range("M132").Value = -33343328
$ws.Range("H136").Value = 3351164.8
$ws.Range("I136").Value = 1667852.6
$ws.Range("J136").Value = 11767726
$ws.Range("K136").Value = 5003557.800000001
$ws.Range("L136").Value = 35303178
$ws.Range("M136").Value = -5001007.800000001
$ws.Range("N136").Value = -35308278

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1453.3462
$ws.Range("I94").Value = 1165.0952
$ws.Range("J94").Value = 2664
$ws.Range("K94").Value = 1165.0952
$ws.Range("L94").Value = 2664
$ws.Range("M94").Value = -714.0952
$ws.Range("N94").Value = -3566

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1455471.2
$ws.Range("I31").Value = 1468.4517
$ws.Range("J31").Value = 5211645
$ws.Range("K31").Value = 1468.4517
$ws.Range("L31").Value = 5211645
$ws.Range("M31").Value = -1173.4517
$ws.Range("N31").Value = -5212235
$ws.Range("H34").Value = 1455471.2
$ws.Range("I34").Value = 1468.4517
$ws.Range("J34").Value = 5211645
$ws.Range("K34").Value = 1468.4517
$ws.Range("L34").Value = 5211645
$ws.Range("M34").Value = -1266.4517
$ws.Range("N34").Value = -5212049
$ws.Range("H58").Value = 1380068.6
$ws.Range("I58").Value = 3157.6924
$ws.Range("J58").Value = 6494309
$ws.Range("K58").Value = 3157.6924
$ws.Range("L58").Value = 6494309
$ws.Range("M58").Value = -2954.6924
$ws.Range("N58").Value = -6494715
$ws.Range("H136").Value = 1380068.6
$ws.Range("I136").Value = 3157.6924
$ws.Range("J136").Value = 6494309
$ws.Range("K136").Value = 9473.0772
$ws.Range("L136").Value = 19482927
$ws.Range("M136").Value = -6923.0772
$ws.Range("N136").Value = -19488027

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 13500
$ws.Range("I80").Value = 5454.5454
$ws.Range("K80").Value = 5454.5454
$ws.Range("M80").Value = -4456.5454
$ws.Range("H83").Value = 13500
$ws.Range("I83").Value = 5454.5454
$ws.Range("K83").Value = 27272.727
$ws.Range("M83").Value = -22280.727
$ws.Range("H97").Value = 22729106
$ws.Range("I97").Value = 2060
$ws.Range("J97").Value = 83334560
$ws.Range("K97").Value = 2060
$ws.Range("L97").Value = 83334560
$ws.Range("M97").Value = -1564
$ws.Range("N97").Value = -83335552
$ws.Range("H126").Value = 13381.25
$ws.Range("I126").Value = 25620
$ws.Range("J126").Value = 1142.5
$ws.Range("K126").Value = 76860
$ws.Range("L126").Value = 3427.5
$ws.Range("M126").Value = -74390
$ws.Range("N126").Value = -8367.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H95").Value = 19695.5
$ws.Range("J95").Value = 19695.5
$ws.Range("L95").Value = 19695.5
$ws.Range("N95").Value = -25187.5
$ws.Range("H96").Value = 31111
$ws.Range("J96").Value = 31111
$ws.Range("L96").Value = 31111
$ws.Range("N96").Value = -36603
$ws.Range("H97").Value = 17736
$ws.Range("J97").Value = 17736
$ws.Range("L97").Value = 17736
$ws.Range("N97").Value = -19718
$ws.Range("H98").Value = 26500
$ws.Range("J98").Value = 26500
$ws.Range("L98").Value = 26500
$ws.Range("N98").Value = -32490
$ws.Range("H99").Value = 35000
$ws.Range("J99").Value = 35000
$ws.Range("L99").Value = 35000
$ws.Range("N99").Value = -40990
$ws.Range("H100").Value = 21786.5
$ws.Range("I100").Value = 32818.688
$ws.Range("J100").Value = 7076.9165
$ws.Range("K100").Value = 32818.688
$ws.Range("L100").Value = 7076.9165
$ws.Range("M100").Value = -32277.688
$ws.Range("N100").Value = -8158.9165
$ws.Range("H101").Value = 20574.857
$ws.Range("J101").Value = 20574.857
$ws.Range("L101").Value = 20574.857
$ws.Range("N101").Value = -27064.857
$ws.Range("H102").Value = 35888
$ws.Range("J102").Value = 35888
$ws.Range("L102").Value = 35888
$ws.Range("N102").Value = -42378
$ws.Range("H103").Value = 32500
$ws.Range("J103").Value = 32500
$ws.Range("L103").Value = 32500
$ws.Range("N103").Value = -34844
$ws.Range("H132").Value = 7145221.5
$ws.Range("I132").Value = 10205139
$ws.Range("J132").Value = 5415.5
$ws.Range("K132").Value = 30615417
$ws.Range("L132").Value = 16246.5
$ws.Range("M132").Value = -30612887
$ws.Range("N132").Value = -21306.5
$ws.Range("H136").Value = 7355953.5
$ws.Range("I136").Value = 14708407
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 44125221
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -44122671
$ws.Range("N136").Value = -15600

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1467.8572
$ws.Range("I96").Value = 1484
$ws.Range("J96").Value = 1333.3334
$ws.Range("K96").Value = 1484
$ws.Range("L96").Value = 1333.3334
$ws.Range("M96").Value = -111
$ws.Range("N96").Value = -4079.3334
$ws.Range("H107").Value = 10783.379
$ws.Range("I107").Value = 12744.1875
$ws.Range("K107").Value = 38232.5625
$ws.Range("M107").Value = -36312.5625
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H132").Value = 1126051.8
$ws.Range("I132").Value = 3912.5881
$ws.Range("J132").Value = 3033688.2
$ws.Range("K132").Value = 11737.7643
$ws.Range("L132").Value = 9101064.600000001
$ws.Range("M132").Value = -9207.764299999999
$ws.Range("N132").Value = -9106124.600000001

Write-Host "Applied 229 cell updates across sheets ALC, ARM, BSM, CRP, GSM, LTW, WVR"
